# Refresh market-price / profit columns (H:N) on each Leve sheet
# with the latest Universalis snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Cells.Item(99, 8).Value = 25642.5
$ws.Cells.Item(99, 9).Value = 50499.5
$ws.Cells.Item(99, 10).Value = 785.5
$ws.Cells.Item(99, 11).Value = 151498.5
$ws.Cells.Item(99, 12).Value = 2356.5
$ws.Cells.Item(99, 13).Value = -150000.5
$ws.Cells.Item(99, 14).Value = -5352.5

# Row 113
$ws.Cells.Item(113, 8).Value = 56908.055
$ws.Cells.Item(113, 9).Value = 92228.82000000001
$ws.Cells.Item(113, 10).Value = 1404
$ws.Cells.Item(113, 11).Value = 92228.82000000001
$ws.Cells.Item(113, 12).Value = 1404
$ws.Cells.Item(113, 13).Value = -88974.82000000001
$ws.Cells.Item(113, 14).Value = -7912

# Row 131
$ws.Cells.Item(131, 8).Value = 4508.7144
$ws.Cells.Item(131, 9).Value = 1173.75
$ws.Cells.Item(131, 10).Value = 4939.032
$ws.Cells.Item(131, 11).Value = 3521.25
$ws.Cells.Item(131, 12).Value = 14817.096
$ws.Cells.Item(131, 13).Value = 1518.75
$ws.Cells.Item(131, 14).Value = -24897.096

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 57389.89
$ws.Cells.Item(2, 9).Value = 1962.5
$ws.Cells.Item(2, 10).Value = 126674.125
$ws.Cells.Item(2, 11).Value = 1962.5
$ws.Cells.Item(2, 12).Value = 126674.125
$ws.Cells.Item(2, 13).Value = -1849.5
$ws.Cells.Item(2, 14).Value = -126900.125

# Row 35
$ws.Cells.Item(35, 8).Value = 837
$ws.Cells.Item(35, 9).Value = 837
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 837
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -431

# Row 61
$ws.Cells.Item(61, 8).Value = 1583.3438
$ws.Cells.Item(61, 9).Value = 1519.4584
$ws.Cells.Item(61, 10).Value = 1775
$ws.Cells.Item(61, 11).Value = 1519.4584
$ws.Cells.Item(61, 12).Value = 1775
$ws.Cells.Item(61, 13).Value = -1307.4584
$ws.Cells.Item(61, 14).Value = -2199

# Row 74
$ws.Cells.Item(74, 8).Value = 7399
$ws.Cells.Item(74, 9).Value = 5000
$ws.Cells.Item(74, 10).Value = 7998.75
$ws.Cells.Item(74, 11).Value = 5000
$ws.Cells.Item(74, 12).Value = 7998.75
$ws.Cells.Item(74, 13).Value = -4126
$ws.Cells.Item(74, 14).Value = -9746.75

# Row 77
$ws.Cells.Item(77, 8).Value = 7399
$ws.Cells.Item(77, 9).Value = 5000
$ws.Cells.Item(77, 10).Value = 7998.75
$ws.Cells.Item(77, 11).Value = 25000
$ws.Cells.Item(77, 12).Value = 39993.75
$ws.Cells.Item(77, 13).Value = -20632
$ws.Cells.Item(77, 14).Value = -48729.75

# Row 116
$ws.Cells.Item(116, 8).Value = 57389.89
$ws.Cells.Item(116, 9).Value = 1962.5
$ws.Cells.Item(116, 10).Value = 126674.125
$ws.Cells.Item(116, 11).Value = 1962.5
$ws.Cells.Item(116, 12).Value = 126674.125
$ws.Cells.Item(116, 13).Value = 331.5
$ws.Cells.Item(116, 14).Value = -131262.125

# Row 132
$ws.Cells.Item(132, 8).Value = 2167.6775
$ws.Cells.Item(132, 9).Value = 1523.3125
$ws.Cells.Item(132, 10).Value = 2855
$ws.Cells.Item(132, 11).Value = 4569.9375
$ws.Cells.Item(132, 12).Value = 8565
$ws.Cells.Item(132, 13).Value = -2039.9375
$ws.Cells.Item(132, 14).Value = -13625

# Row 136
$ws.Cells.Item(136, 8).Value = 1583.3438
$ws.Cells.Item(136, 9).Value = 1519.4584
$ws.Cells.Item(136, 10).Value = 1775
$ws.Cells.Item(136, 11).Value = 4558.3752
$ws.Cells.Item(136, 12).Value = 5325
$ws.Cells.Item(136, 13).Value = -2008.3752
$ws.Cells.Item(136, 14).Value = -10425

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 57389.89
$ws.Cells.Item(3, 9).Value = 1962.5
$ws.Cells.Item(3, 10).Value = 126674.125
$ws.Cells.Item(3, 11).Value = 1962.5
$ws.Cells.Item(3, 12).Value = 126674.125
$ws.Cells.Item(3, 13).Value = -1848.5
$ws.Cells.Item(3, 14).Value = -126902.125

# Row 36
$ws.Cells.Item(36, 8).Value = 13562.2
$ws.Cells.Item(36, 9).Value = 13562.2
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 13562.2
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -13028.2
$ws.Cells.Item(36, 14).Value = $null

# Row 134
$ws.Cells.Item(134, 8).Value = 2622.8333
$ws.Cells.Item(134, 9).Value = 2664.4644
$ws.Cells.Item(134, 10).Value = 2040
$ws.Cells.Item(134, 11).Value = 7993.3932
$ws.Cells.Item(134, 12).Value = 6120
$ws.Cells.Item(134, 13).Value = -5458.3932
$ws.Cells.Item(134, 14).Value = -11190

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Cells.Item(99, 8).Value = 24798.8
$ws.Cells.Item(99, 9).Value = 6326.6665
$ws.Cells.Item(99, 10).Value = 52507
$ws.Cells.Item(99, 11).Value = 6326.6665
$ws.Cells.Item(99, 12).Value = 52507
$ws.Cells.Item(99, 13).Value = -4828.6665
$ws.Cells.Item(99, 14).Value = -55503

# Row 107
$ws.Cells.Item(107, 8).Value = 1238.8462
$ws.Cells.Item(107, 9).Value = 1772.1666
$ws.Cells.Item(107, 10).Value = 781.7143
$ws.Cells.Item(107, 11).Value = 1772.1666
$ws.Cells.Item(107, 12).Value = 781.7143
$ws.Cells.Item(107, 13).Value = 147.8334
$ws.Cells.Item(107, 14).Value = -4621.7143

# Row 122
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = $null
$ws.Cells.Item(122, 14).Value = $null

# Row 126
$ws.Cells.Item(126, 8).Value = 24798.8
$ws.Cells.Item(126, 9).Value = 6326.6665
$ws.Cells.Item(126, 10).Value = 52507
$ws.Cells.Item(126, 11).Value = 18979.9995
$ws.Cells.Item(126, 12).Value = 157521
$ws.Cells.Item(126, 13).Value = -16509.9995
$ws.Cells.Item(126, 14).Value = -162461

# Row 132
$ws.Cells.Item(132, 8).Value = 4336.8237
$ws.Cells.Item(132, 9).Value = 4279.0835
$ws.Cells.Item(132, 10).Value = 4475.4
$ws.Cells.Item(132, 11).Value = 12837.2505
$ws.Cells.Item(132, 12).Value = 13426.2
$ws.Cells.Item(132, 13).Value = -10307.2505
$ws.Cells.Item(132, 14).Value = -18486.2

# Row 134
$ws.Cells.Item(134, 8).Value = 1166.4849
$ws.Cells.Item(134, 9).Value = 1086
$ws.Cells.Item(134, 10).Value = 1750
$ws.Cells.Item(134, 11).Value = 3258
$ws.Cells.Item(134, 12).Value = 5250
$ws.Cells.Item(134, 13).Value = -723
$ws.Cells.Item(134, 14).Value = -10320

# Row 138
$ws.Cells.Item(138, 8).Value = 62092.5
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 62092.5
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 62092.5
$ws.Cells.Item(138, 14).Value = -72372.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 7954.143
$ws.Cells.Item(5, 9).Value = 727.2105
$ws.Cells.Item(5, 10).Value = 23211
$ws.Cells.Item(5, 11).Value = 2181.6315
$ws.Cells.Item(5, 12).Value = 69633
$ws.Cells.Item(5, 13).Value = -2069.6315
$ws.Cells.Item(5, 14).Value = -69857

# Row 34
$ws.Cells.Item(34, 8).Value = 504
$ws.Cells.Item(34, 9).Value = 130
$ws.Cells.Item(34, 10).Value = 2000
$ws.Cells.Item(34, 11).Value = 390
$ws.Cells.Item(34, 12).Value = 6000
$ws.Cells.Item(34, 13).Value = -306
$ws.Cells.Item(34, 14).Value = -6168

# Row 80
$ws.Cells.Item(80, 8).Value = 100003
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 100003
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 300009
$ws.Cells.Item(80, 14).Value = -301881

# Row 83
$ws.Cells.Item(83, 8).Value = 100003
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 100003
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 900027
$ws.Cells.Item(83, 14).Value = -909387

# Row 122
$ws.Cells.Item(122, 8).Value = 6331.647
$ws.Cells.Item(122, 9).Value = 376
$ws.Cells.Item(122, 10).Value = 50999
$ws.Cells.Item(122, 11).Value = 3384
$ws.Cells.Item(122, 12).Value = 458991
$ws.Cells.Item(122, 13).Value = -934
$ws.Cells.Item(122, 14).Value = -463891

# Row 135
$ws.Cells.Item(135, 8).Value = 7954.143
$ws.Cells.Item(135, 9).Value = 727.2105
$ws.Cells.Item(135, 10).Value = 23211
$ws.Cells.Item(135, 11).Value = 6544.8945
$ws.Cells.Item(135, 12).Value = 208899
$ws.Cells.Item(135, 13).Value = -4009.8945
$ws.Cells.Item(135, 14).Value = -213969

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Cells.Item(97, 8).Value = 32259442
$ws.Cells.Item(97, 9).Value = 33334658
$ws.Cells.Item(97, 10).Value = 3000
$ws.Cells.Item(97, 11).Value = 33334658
$ws.Cells.Item(97, 12).Value = 3000
$ws.Cells.Item(97, 13).Value = -33334162
$ws.Cells.Item(97, 14).Value = -3992

# Row 102
$ws.Cells.Item(102, 8).Value = 2189
$ws.Cells.Item(102, 9).Value = 2161.2104
$ws.Cells.Item(102, 10).Value = 2247.6667
$ws.Cells.Item(102, 11).Value = 2161.2104
$ws.Cells.Item(102, 12).Value = 2247.6667
$ws.Cells.Item(102, 13).Value = -539.2103999999999
$ws.Cells.Item(102, 14).Value = -5491.6667

# Row 132
$ws.Cells.Item(132, 8).Value = 4874.75
$ws.Cells.Item(132, 9).Value = 4666.6665
$ws.Cells.Item(132, 10).Value = 4999.6
$ws.Cells.Item(132, 11).Value = 13999.9995
$ws.Cells.Item(132, 12).Value = 14998.8
$ws.Cells.Item(132, 13).Value = -11469.9995
$ws.Cells.Item(132, 14).Value = -20058.8

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Cells.Item(10, 8).Value = 5000
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 5000
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 5000
$ws.Cells.Item(10, 13).Value = $null
$ws.Cells.Item(10, 14).Value = -5280

# Row 40
$ws.Cells.Item(40, 8).Value = 2374.75
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 2374.75
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 2374.75
$ws.Cells.Item(40, 13).Value = $null
$ws.Cells.Item(40, 14).Value = -2646.75

# Row 132
$ws.Cells.Item(132, 8).Value = 4473.2354
$ws.Cells.Item(132, 9).Value = 4536.067
$ws.Cells.Item(132, 10).Value = 4002
$ws.Cells.Item(132, 11).Value = 13608.201
$ws.Cells.Item(132, 12).Value = 12006
$ws.Cells.Item(132, 13).Value = -11078.201
$ws.Cells.Item(132, 14).Value = -17066

# Row 133
$ws.Cells.Item(133, 8).Value = 59895
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 59895
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 59895
$ws.Cells.Item(133, 14).Value = -64955

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 2188.3635
$ws.Cells.Item(122, 9).Value = 1302
$ws.Cells.Item(122, 10).Value = 2385.3333
$ws.Cells.Item(122, 11).Value = 3906
$ws.Cells.Item(122, 12).Value = 7155.999899999999
$ws.Cells.Item(122, 13).Value = -1456
$ws.Cells.Item(122, 14).Value = -12055.9999

# Row 132
$ws.Cells.Item(132, 8).Value = 2333
$ws.Cells.Item(132, 9).Value = 1000
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 3000
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -470
$ws.Cells.Item(132, 14).Value = -20057

# Row 138
$ws.Cells.Item(138, 8).Value = 55994.285
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 55994.285
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 55994.285
$ws.Cells.Item(138, 14).Value = -66274.285
